$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (mean), C (CI_low), D (CI_high) for rows 2-14
$values = @{
    2  = @(0.013422399255276, 0.000608678669487831, 0.0262361198410642)
    3  = @(0.0183425262321837, -0.0000472323760656104, 0.0367322848404331)
    4  = @(0.00303082617256612, -0.041811584630556, 0.0478732369756882)
    5  = @(0.0393580937975179, 0.00268962283031415, 0.0760265647647216)
    6  = @(0.0273693272733676, -0.0170822407580437, 0.0718208953047788)
    7  = @(0.00499455885299538, -0.0522255970791038, 0.0622147147850945)
    8  = @(0.0584586119008044, 0.00203979381560638, 0.114877429986002)
    9  = @(-0.00232934159943271, -0.0470463441134809, 0.0423876609146155)
    10 = @(-0.0113963386871638, -0.0702712806554948, 0.0474786032811672)
    11 = @(0.0131981795918118, -0.0127251092327356, 0.0391214684163592)
    12 = @(-0.0141821128475735, -0.0579220523452886, 0.0295578266501417)
    13 = @(-0.023636339470866, -0.0721256916264213, 0.0248530126846894)
    14 = @(0.0199605289964537, -0.0039850405411648, 0.0439060985340722)
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 2).Value = $v[0]
    $ws.Cells.Item($row, 3).Value = $v[1]
    $ws.Cells.Item($row, 4).Value = $v[2]
}
